$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B (associated_text) values for rows 2-8
$ws.Range("B2").Value = "for"
$ws.Range("B3").Value = "for"
$ws.Range("B4").Value = "FOR"
$ws.Range("B5").Value = "for"
$ws.Range("B6").Value = "for"
$ws.Range("B7").Value = "for"
$ws.Range("B8").Value = "for"

# Swap the URLs in A7 and A8
$a7 = $ws.Range("A7").Value2
$a8 = $ws.Range("A8").Value2
$ws.Range("A7").Value = $a8
$ws.Range("A8").Value = $a7
